$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 3477.2
$ws.Range("I11").Value = 3477.2
$ws.Range("K11").Value = 3477.2
$ws.Range("M11").Value = -3337.2
$ws.Range("H28").Value = 794.70966
$ws.Range("I28").Value = 780.5789
$ws.Range("J28").Value = 817.0833
$ws.Range("K28").Value = 780.5789
$ws.Range("L28").Value = 817.0833
$ws.Range("M28").Value = -295.5789
$ws.Range("N28").Value = -1787.0833
$ws.Range("H105").Value = 16671
$ws.Range("J105").Value = 16671
$ws.Range("L105").Value = 16671
$ws.Range("N105").Value = -23659
$ws.Range("H132").Value = 405095.44
$ws.Range("I132").Value = 421944.4
$ws.Range("J132").Value = 720
$ws.Range("K132").Value = 1265833.2
$ws.Range("L132").Value = 2160
$ws.Range("M132").Value = -1263303.2
$ws.Range("N132").Value = -7220
$ws.Range("H137").Value = 21740404
$ws.Range("I137").Value = 1149
$ws.Range("J137").Value = 76924670
$ws.Range("K137").Value = 3447
$ws.Range("L137").Value = 230774010
$ws.Range("M137").Value = -897
$ws.Range("N137").Value = -230779110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1401839.6
$ws.Range("I2").Value = 1186.9
$ws.Range("J2").Value = 2675160.2
$ws.Range("K2").Value = 1186.9
$ws.Range("L2").Value = 2675160.2
$ws.Range("M2").Value = -1073.9
$ws.Range("N2").Value = -2675386.2
$ws.Range("H32").Value = 3864.0645
$ws.Range("I32").Value = 3871.6829
$ws.Range("J32").Value = 3807.2727
$ws.Range("K32").Value = 3871.6829
$ws.Range("L32").Value = 3807.2727
$ws.Range("M32").Value = -3584.6829
$ws.Range("N32").Value = -4381.2727
$ws.Range("H74").Value = 14710676
$ws.Range("I74").Value = 19231600
$ws.Range("J74").Value = 17673.625
$ws.Range("K74").Value = 19231600
$ws.Range("L74").Value = 17673.625
$ws.Range("M74").Value = -19230726
$ws.Range("N74").Value = -19421.625
$ws.Range("H77").Value = 14710676
$ws.Range("I77").Value = 19231600
$ws.Range("J77").Value = 17673.625
$ws.Range("K77").Value = 96158000
$ws.Range("L77").Value = 88368.125
$ws.Range("M77").Value = -96153632
$ws.Range("N77").Value = -97104.125
$ws.Range("H116").Value = 1401839.6
$ws.Range("I116").Value = 1186.9
$ws.Range("J116").Value = 2675160.2
$ws.Range("K116").Value = 1186.9
$ws.Range("L116").Value = 2675160.2
$ws.Range("M116").Value = 1107.1
$ws.Range("N116").Value = -2679748.2
$ws.Range("H132").Value = 790332.5
$ws.Range("I132").Value = 1114294.1
$ws.Range("J132").Value = 72989.07000000001
$ws.Range("K132").Value = 3342882.3
$ws.Range("L132").Value = 218967.21
$ws.Range("M132").Value = -3340352.3
$ws.Range("N132").Value = -224027.21
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1401839.6
$ws.Range("I3").Value = 1186.9
$ws.Range("J3").Value = 2675160.2
$ws.Range("K3").Value = 1186.9
$ws.Range("L3").Value = 2675160.2
$ws.Range("M3").Value = -1072.9
$ws.Range("N3").Value = -2675388.2
$ws.Range("H86").Value = 2082.54
$ws.Range("I86").Value = 1865.2941
$ws.Range("J86").Value = 2544.1875
$ws.Range("K86").Value = 1865.2941
$ws.Range("L86").Value = 2544.1875
$ws.Range("M86").Value = -742.2941000000001
$ws.Range("N86").Value = -4790.1875
$ws.Range("H89").Value = 2082.54
$ws.Range("I89").Value = 1865.2941
$ws.Range("J89").Value = 2544.1875
$ws.Range("K89").Value = 9326.470499999999
$ws.Range("L89").Value = 12720.9375
$ws.Range("M89").Value = -3710.470499999999
$ws.Range("N89").Value = -23952.9375
$ws.Range("H134").Value = 3726605.5
$ws.Range("I134").Value = 4471708.5
$ws.Range("J134").Value = 1090.3334
$ws.Range("K134").Value = 13415125.5
$ws.Range("L134").Value = 3271.0002
$ws.Range("M134").Value = -13412590.5
$ws.Range("N134").Value = -8341.0002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1980.7142
$ws.Range("I132").Value = 1918.2778
$ws.Range("K132").Value = 5754.8334
$ws.Range("M132").Value = -3224.8334
$ws.Range("H134").Value = 1464.973
$ws.Range("I134").Value = 1688.2222
$ws.Range("J134").Value = 862.2
$ws.Range("K134").Value = 5064.6666
$ws.Range("L134").Value = 2586.6
$ws.Range("M134").Value = -2529.6666
$ws.Range("N134").Value = -7656.6
$ws.Range("H141").Value = 35492.35
$ws.Range("J141").Value = 35492.35
$ws.Range("L141").Value = 35492.35
$ws.Range("N141").Value = -45852.35
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1615.5
$ws.Range("I97").Value = 1594.375
$ws.Range("J97").Value = 1700
$ws.Range("K97").Value = 1594.375
$ws.Range("L97").Value = 1700
$ws.Range("M97").Value = -1098.375
$ws.Range("N97").Value = -2692
$ws.Range("H107").Value = 773.4
$ws.Range("I107").Value = 508.66666
$ws.Range("J107").Value = 949.8889
$ws.Range("K107").Value = 508.66666
$ws.Range("L107").Value = 949.8889
$ws.Range("M107").Value = 1411.33334
$ws.Range("N107").Value = -4789.8889
$ws.Range("H134").Value = 28775
$ws.Range("J134").Value = 28775
$ws.Range("L134").Value = 86325
$ws.Range("N134").Value = -91395
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1361.2778
$ws.Range("I46").Value = 925.0833
$ws.Range("J46").Value = 2233.6667
$ws.Range("K46").Value = 925.0833
$ws.Range("L46").Value = 2233.6667
$ws.Range("M46").Value = -737.0833
$ws.Range("N46").Value = -2609.6667
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H136").Value = 5375.433
$ws.Range("I136").Value = 6380.15
$ws.Range("J136").Value = 3366
$ws.Range("K136").Value = 19140.45
$ws.Range("L136").Value = 10098
$ws.Range("M136").Value = -16590.45
$ws.Range("N136").Value = -15198
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3555.157
$ws.Range("I132").Value = 4166.1904
$ws.Range("J132").Value = 703.6667
$ws.Range("K132").Value = 12498.5712
$ws.Range("L132").Value = 2111.0001
$ws.Range("M132").Value = -9968.571200000002
$ws.Range("N132").Value = -7171.0001
$ws.Range("H136").Value = 6340.343
$ws.Range("I136").Value = 7030.0967
$ws.Range("J136").Value = 994.75
$ws.Range("K136").Value = 21090.2901
$ws.Range("L136").Value = 2984.25
$ws.Range("M136").Value = -18540.2901
$ws.Range("N136").Value = -8084.25

Write-Output "Applied all cell updates"